$d = $word.ActiveDocument

$d.Content.Find.Execute("Operating System: Linux 6.14.0-1011-aws", $true, $false, $false, $false, $false, $true, 1, $false, "Operating System: Linux 6.11.0-1018-azure", 2)
$d.Content.Find.Execute("Python Version: 3.12.3", $true, $false, $false, $false, $false, $true, 1, $false, "Python Version: 3.12.12", 2)
$d.Content.Find.Execute("Total Memory (MB): 7820.98", $true, $false, $false, $false, $false, $true, 1, $false, "Total Memory (MB): 15995.58", 2)
$d.Content.Find.Execute("Available Memory (MB): 7008.12", $true, $false, $false, $false, $false, $true, 1, $false, "Available Memory (MB): 14668.07", 2)
